# Sample Project / Main.xlsx — "Rules" sheet edit
#
# Rule row 11 (B11:E11) had its Rule-ID changed from the text "R40" to the
# text "1". Like the other Rule-ID cells in column B (R10, R20, R30, ...)
# this is a text label, not a number, so it must be written as text even
# though "1" looks numeric — otherwise Excel's automatic type detection
# would store it as a number instead of a shared string.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Leading apostrophe forces Excel to keep the entry as text ("1") rather
# than auto-converting it to the number 1.
$ws.Range("B11").Value = "'1"
